$wb = $excel.ActiveWorkbook

# --- Sheet "Calificaciones": fill in previously-missing 2P/Final grades for TEMAS DE ADMINISTRACION (col L / col X) ---
$wsCal = $wb.Worksheets.Item("Calificaciones")
$wsCal.Range("L7").Value = 10
$wsCal.Range("X7").Value = 9
$wsCal.Range("L12").Value = 9
$wsCal.Range("X12").Value = 9
$wsCal.Range("L15").Value = 9
$wsCal.Range("L18").Value = 8
$wsCal.Range("L19").Value = 9
$wsCal.Range("L20").Value = 9
$wsCal.Range("L22").Value = 8
$wsCal.Range("L23").Value = 8
$wsCal.Range("X23").Value = 8
$wsCal.Range("L24").Value = 8
$wsCal.Range("X24").Value = 8
$wsCal.Range("L28").Value = 8
$wsCal.Range("X28").Value = 9
$wsCal.Range("L29").Value = 8
$wsCal.Range("X29").Value = 8
$wsCal.Range("L30").Value = 8
$wsCal.Range("X30").Value = 7
$wsCal.Range("L31").Value = 8
$wsCal.Range("L32").Value = 8
$wsCal.Range("X32").Value = 7
$wsCal.Range("L35").Value = 8
$wsCal.Range("L36").Value = 10
$wsCal.Range("L42").Value = 10
$wsCal.Range("X42").Value = 9

# --- Sheet "Asistencias": recompute F2/F3 attendance percentage for TEMAS DE ADMINISTRACION (col L / col R) ---
$wsAsi = $wb.Worksheets.Item("Asistencias")
$wsAsi.Range("L7").Value = 81.59999999999999
$wsAsi.Range("R7").Value = 81.59999999999999
$wsAsi.Range("L12").Value = 100
$wsAsi.Range("R12").Value = 100
$wsAsi.Range("L18").Value = 83.7
$wsAsi.Range("R18").Value = 83.7
$wsAsi.Range("L19").Value = 89.8
$wsAsi.Range("R19").Value = 89.8
$wsAsi.Range("L20").Value = 81.59999999999999
$wsAsi.Range("R20").Value = 81.59999999999999
$wsAsi.Range("L22").Value = 87.8
$wsAsi.Range("R22").Value = 87.8
$wsAsi.Range("L23").Value = 81.59999999999999
$wsAsi.Range("R23").Value = 81.59999999999999
$wsAsi.Range("L29").Value = 81.59999999999999
$wsAsi.Range("R29").Value = 81.59999999999999
$wsAsi.Range("L30").Value = 81.59999999999999
$wsAsi.Range("R30").Value = 81.59999999999999
$wsAsi.Range("L31").Value = 81.59999999999999
$wsAsi.Range("R31").Value = 81.59999999999999
$wsAsi.Range("L32").Value = 81.59999999999999
$wsAsi.Range("R32").Value = 81.59999999999999
$wsAsi.Range("L35").Value = 81.59999999999999
$wsAsi.Range("R35").Value = 81.59999999999999
$wsAsi.Range("L42").Value = 95.90000000000001
$wsAsi.Range("R42").Value = 95.90000000000001

# --- Sheet "Totales": the three summary rows (TEMAS DE ADMINISTRACION / INTRODUCCION A LA ECONOMIA /
#     TEMAS DE FILOSOFIA) are re-ordered and their stats recomputed now that 2P grades are in ---
$wsTot = $wb.Worksheets.Item("Totales")

$wsTot.Range("A4").Value = "INTRODUCCIÓN A LA ECONOMÍA"
$wsTot.Range("B4").Value = "Miguel Cruz Nayeli Nayreth"
$wsTot.Range("H4").Value = 8.199999999999999

$wsTot.Range("A5").Value = "TEMAS DE FILOSOFÍA"
$wsTot.Range("B5").Value = "Barrientos Ortiz Yuliana Isabel"
$wsTot.Range("D5").Value = 39
$wsTot.Range("E5").Value = 0
$wsTot.Range("F5").Value = 100
$wsTot.Range("G5").Value = 0
$wsTot.Range("H5").Value = 9

$wsTot.Range("A6").Value = "TEMAS DE ADMINISTRACIÓN"
$wsTot.Range("B6").Value = "Saucedo Rivalcoba Liliana Guadalupe"
$wsTot.Range("H6").Value = 9.300000000000001

# --- Sheet "Rescatables": the two students who were pending TEMAS DE ADMINISTRACION (2P) now
#     passed it, so their "rescatable" rows are removed ---
$wsRes = $wb.Worksheets.Item("Rescatables")
$wsRes.Rows.Item(9).Delete()
$wsRes.Rows.Item(5).Delete()
